$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Preserve the (visually-identical) highlighted style that currently lives on
# C46 by copying it onto the new C47 cell before C46's own style is toned
# down below.
$ws.Range("C46").Copy($ws.Range("C47"))

# B46/C46 lose their (no-op) fill flag - they pick up the plain styles that
# are already used elsewhere in the B/C columns (e.g. row 45).
$ws.Range("B45").Copy()
$ws.Range("B46").PasteSpecial(-4122)

$ws.Range("C45").Copy()
$ws.Range("C46").PasteSpecial(-4122)

# New row 47: date/hours formatting copied from row 46 (A47, B47).
$ws.Range("A46").Copy()
$ws.Range("A47").PasteSpecial(-4122)

$ws.Range("B45").Copy()
$ws.Range("B47").PasteSpecial(-4122)

# Values for the new work-log entry.
$ws.Range("A47").Value = 45751
$ws.Range("B47").Value = 4
$ws.Range("C47").Value = "Switched Ollama to higher version and modified the codes to minimize scanning time"

# Update the saved view state to match what Excel recorded after the edit.
$ws.Range("B1").Select()
$ws.Application.ActiveWindow.ScrollRow = 1
$ws.Range("C49").Select()
